$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Tester 2 name" column is being inserted at G (pushing the existing
# Tester-2 Pass/Fail result from G into H). For every data row (4-10):
#   1. Copy the current G cell (value + conditional Pass/Fail fill) into H.
#   2. Copy F's plain (no-fill) formatting into G, then overwrite G's value
#      with the new tester name so G keeps the unstyled look.
$rows = 4..10
foreach ($r in $rows) {
    $gCell = $ws.Range("G$r")
    $hCell = $ws.Range("H$r")
    $fCell = $ws.Range("F$r")

    $gCell.Copy($hCell)
    $fCell.Copy($gCell)
    $gCell.Value = "Aneesh Dalvi"
}

# Reset the saved scroll position so the sheet view no longer pins topLeftCell=A4.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
